$d = $word.ActiveDocument

# --- Change 1: "Pe" + hidden bookmark + "rmohonan" -> "Permohonan" ---
# (merges the two split runs back into a single word and, by relocating the
# _GoBack bookmark below, removes the stray bookmark pair here)
$d.Content.Find.Execute("Pe" + [char]0x0008 + "rmohonan", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0) | Out-Null

# --- Change 2: strike-through the whole "Dashboard Penyelaras : pada ruang ..." bullet ---
# Locate the paragraph by its distinctive lead-in text.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Dashboard Penyelaras : pada ruang")) {
        $target = $p
        break
    }
}

$full = $target.Range
$full.Font.StrikeThrough = 1

# Move the (hidden) _GoBack bookmark so it wraps this whole paragraph's text,
# which removes it from its old location next to "Pe"/"rmohonan" above.
$bkRange = $d.Range($target.Range.Start, $target.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $bkRange) | Out-Null
